$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: 5.5 -> 5.6
$ws.Range("B2").Value = 5.6

# Update C4: 1.2 -> 1.25
$ws.Range("C4").Value = 1.25

# Update the active selection to B2
$ws.Range("B2").Select()
